$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Re-shuffle match data (columns F:V) among rows that share the
#    same data_partida timestamp. Columns A:E are left untouched.
#    Done with live Range.Value2 reads so the written values always
#    come straight from the workbook itself.
# ---------------------------------------------------------------

# group: 78, 80, 81
$row78 = $ws.Range("F78:V78").Value2
$row80 = $ws.Range("F80:V80").Value2
$row81 = $ws.Range("F81:V81").Value2
$ws.Range("F78:V78").Value2 = $row81
$ws.Range("F80:V80").Value2 = $row78
$ws.Range("F81:V81").Value2 = $row80

# group: 192, 193
$row192 = $ws.Range("F192:V192").Value2
$row193 = $ws.Range("F193:V193").Value2
$ws.Range("F192:V192").Value2 = $row193
$ws.Range("F193:V193").Value2 = $row192

# group: 198, 199
$row198 = $ws.Range("F198:V198").Value2
$row199 = $ws.Range("F199:V199").Value2
$ws.Range("F198:V198").Value2 = $row199
$ws.Range("F199:V199").Value2 = $row198

# group: 203, 205
$row203 = $ws.Range("F203:V203").Value2
$row205 = $ws.Range("F205:V205").Value2
$ws.Range("F203:V203").Value2 = $row205
$ws.Range("F205:V205").Value2 = $row203

# group: 206, 207, 208, 209
$row206 = $ws.Range("F206:V206").Value2
$row207 = $ws.Range("F207:V207").Value2
$row208 = $ws.Range("F208:V208").Value2
$row209 = $ws.Range("F209:V209").Value2
$ws.Range("F206:V206").Value2 = $row209
$ws.Range("F207:V207").Value2 = $row208
$ws.Range("F208:V208").Value2 = $row206
$ws.Range("F209:V209").Value2 = $row207

# group: 214, 215
$row214 = $ws.Range("F214:V214").Value2
$row215 = $ws.Range("F215:V215").Value2
$ws.Range("F214:V214").Value2 = $row215
$ws.Range("F215:V215").Value2 = $row214

# ---------------------------------------------------------------
# 2) Append the 8 new matches scraped on 31/10/2023 as rows 218..225
#    (Indice 217..224). Styles for columns A (bold/border) and E
#    (datetime format) are copied from the last existing row (217)
#    so the new rows match the sheet formatting exactly.
# ---------------------------------------------------------------

$ws.Range("A217:E217").Copy()
$ws.Range("A218:E225").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(218, 1).Value2 = 217
$ws.Cells.Item(218, 2).Value2 = "china"
$ws.Cells.Item(218, 3).Value2 = "jia-league"
$ws.Cells.Item(218, 4).Value2 = "2023"
$ws.Cells.Item(218, 5).Value2 = 45227.35416666666
$ws.Cells.Item(218, 6).Value2 = "Suzhou Dongwu"
$ws.Cells.Item(218, 7).Value2 = 1
$ws.Cells.Item(218, 8).Value2 = "Guangxi Pingguo Haliao"
$ws.Cells.Item(218, 9).Value2 = 3
$ws.Cells.Item(218, 10).Value2 = 4
$ws.Cells.Item(218, 11).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(218, 12).Value2 = 4.2
$ws.Cells.Item(218, 13).Value2 = "28/10/2023 08:16"
$ws.Cells.Item(218, 14).Value2 = 3.08
$ws.Cells.Item(218, 15).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(218, 16).Value2 = 3.12
$ws.Cells.Item(218, 17).Value2 = "28/10/2023 08:16"
$ws.Cells.Item(218, 18).Value2 = 1.9
$ws.Cells.Item(218, 19).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(218, 20).Value2 = 2
$ws.Cells.Item(218, 21).Value2 = "28/10/2023 08:16"
$ws.Cells.Item(218, 22).Value2 = "https://www.betexplorer.com/football/china/jia-league/suzhou-dongwu-guangxi-pingguo-haliao/EoR6mzpO/"

$ws.Cells.Item(219, 1).Value2 = 218
$ws.Cells.Item(219, 2).Value2 = "china"
$ws.Cells.Item(219, 3).Value2 = "jia-league"
$ws.Cells.Item(219, 4).Value2 = "2023"
$ws.Cells.Item(219, 5).Value2 = 45227.35416666666
$ws.Cells.Item(219, 6).Value2 = "Dandong Tengyue"
$ws.Cells.Item(219, 7).Value2 = 1
$ws.Cells.Item(219, 8).Value2 = "Shijiazhuang Gongfu"
$ws.Cells.Item(219, 9).Value2 = 2
$ws.Cells.Item(219, 10).Value2 = 3.93
$ws.Cells.Item(219, 11).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(219, 12).Value2 = 8.14
$ws.Cells.Item(219, 13).Value2 = "28/10/2023 08:29"
$ws.Cells.Item(219, 14).Value2 = 3.13
$ws.Cells.Item(219, 15).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(219, 16).Value2 = 5.09
$ws.Cells.Item(219, 17).Value2 = "28/10/2023 08:29"
$ws.Cells.Item(219, 18).Value2 = 1.89
$ws.Cells.Item(219, 19).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(219, 20).Value2 = 1.36
$ws.Cells.Item(219, 21).Value2 = "28/10/2023 08:28"
$ws.Cells.Item(219, 22).Value2 = "https://www.betexplorer.com/football/china/jia-league/dandong-tengyue-shijiazhuang-gongfu/lb8S4yFb/"

$ws.Cells.Item(220, 1).Value2 = 219
$ws.Cells.Item(220, 2).Value2 = "china"
$ws.Cells.Item(220, 3).Value2 = "jia-league"
$ws.Cells.Item(220, 4).Value2 = "2023"
$ws.Cells.Item(220, 5).Value2 = 45227.35416666666
$ws.Cells.Item(220, 6).Value2 = "Guangzhou FC"
$ws.Cells.Item(220, 7).Value2 = 0
$ws.Cells.Item(220, 8).Value2 = "Qingdao West Coast"
$ws.Cells.Item(220, 9).Value2 = 2
$ws.Cells.Item(220, 10).Value2 = 4.59
$ws.Cells.Item(220, 11).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(220, 12).Value2 = 5.21
$ws.Cells.Item(220, 13).Value2 = "28/10/2023 08:29"
$ws.Cells.Item(220, 14).Value2 = 3.8
$ws.Cells.Item(220, 15).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(220, 16).Value2 = 4.47
$ws.Cells.Item(220, 17).Value2 = "28/10/2023 08:29"
$ws.Cells.Item(220, 18).Value2 = 1.64
$ws.Cells.Item(220, 19).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(220, 20).Value2 = 1.56
$ws.Cells.Item(220, 21).Value2 = "28/10/2023 08:29"
$ws.Cells.Item(220, 22).Value2 = "https://www.betexplorer.com/football/china/jia-league/guangzhou-fc-qingdao-west-coast/pCM2lGUH/"

$ws.Cells.Item(221, 1).Value2 = 220
$ws.Cells.Item(221, 2).Value2 = "china"
$ws.Cells.Item(221, 3).Value2 = "jia-league"
$ws.Cells.Item(221, 4).Value2 = "2023"
$ws.Cells.Item(221, 5).Value2 = 45227.35416666666
$ws.Cells.Item(221, 6).Value2 = "Heilongjiang Ice City"
$ws.Cells.Item(221, 7).Value2 = 1
$ws.Cells.Item(221, 8).Value2 = "Dongguan Guanlian"
$ws.Cells.Item(221, 9).Value2 = 0
$ws.Cells.Item(221, 10).Value2 = 1.98
$ws.Cells.Item(221, 11).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(221, 12).Value2 = 1.31
$ws.Cells.Item(221, 13).Value2 = "28/10/2023 08:28"
$ws.Cells.Item(221, 14).Value2 = 3.19
$ws.Cells.Item(221, 15).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(221, 16).Value2 = 5.37
$ws.Cells.Item(221, 17).Value2 = "28/10/2023 08:28"
$ws.Cells.Item(221, 18).Value2 = 3.54
$ws.Cells.Item(221, 19).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(221, 20).Value2 = 9.66
$ws.Cells.Item(221, 21).Value2 = "28/10/2023 08:28"
$ws.Cells.Item(221, 22).Value2 = "https://www.betexplorer.com/football/china/jia-league/heilongjiang-ice-city-dongguan-guanlian/IaiH7wVu/"

$ws.Cells.Item(222, 1).Value2 = 221
$ws.Cells.Item(222, 2).Value2 = "china"
$ws.Cells.Item(222, 3).Value2 = "jia-league"
$ws.Cells.Item(222, 4).Value2 = "2023"
$ws.Cells.Item(222, 5).Value2 = 45227.35416666666
$ws.Cells.Item(222, 6).Value2 = "Jiangxi Lushan"
$ws.Cells.Item(222, 7).Value2 = 3
$ws.Cells.Item(222, 8).Value2 = "Shenyang Urban FC"
$ws.Cells.Item(222, 9).Value2 = 1
$ws.Cells.Item(222, 10).Value2 = 2.96
$ws.Cells.Item(222, 11).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(222, 12).Value2 = 1.98
$ws.Cells.Item(222, 13).Value2 = "28/10/2023 08:04"
$ws.Cells.Item(222, 14).Value2 = 3.19
$ws.Cells.Item(222, 15).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(222, 16).Value2 = 3.58
$ws.Cells.Item(222, 17).Value2 = "28/10/2023 08:04"
$ws.Cells.Item(222, 18).Value2 = 2.27
$ws.Cells.Item(222, 19).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(222, 20).Value2 = 3.64
$ws.Cells.Item(222, 21).Value2 = "28/10/2023 08:04"
$ws.Cells.Item(222, 22).Value2 = "https://www.betexplorer.com/football/china/jia-league/jiangxi-lushan-shenyang-urban-fc/GpAK6cpn/"

$ws.Cells.Item(223, 1).Value2 = 222
$ws.Cells.Item(223, 2).Value2 = "china"
$ws.Cells.Item(223, 3).Value2 = "jia-league"
$ws.Cells.Item(223, 4).Value2 = "2023"
$ws.Cells.Item(223, 5).Value2 = 45227.35416666666
$ws.Cells.Item(223, 6).Value2 = "Jinan Xingzhou"
$ws.Cells.Item(223, 7).Value2 = 1
$ws.Cells.Item(223, 8).Value2 = "Wuxi Wugou"
$ws.Cells.Item(223, 9).Value2 = 2
$ws.Cells.Item(223, 10).Value2 = 1.44
$ws.Cells.Item(223, 11).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(223, 12).Value2 = 1.5
$ws.Cells.Item(223, 13).Value2 = "28/10/2023 08:02"
$ws.Cells.Item(223, 14).Value2 = 4.06
$ws.Cells.Item(223, 15).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(223, 16).Value2 = 4.3
$ws.Cells.Item(223, 17).Value2 = "28/10/2023 08:02"
$ws.Cells.Item(223, 18).Value2 = 6.18
$ws.Cells.Item(223, 19).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(223, 20).Value2 = 6.39
$ws.Cells.Item(223, 21).Value2 = "28/10/2023 08:02"
$ws.Cells.Item(223, 22).Value2 = "https://www.betexplorer.com/football/china/jia-league/jinan-xingzhou-wuxi-wugou/A79O5Hah/"

$ws.Cells.Item(224, 1).Value2 = 223
$ws.Cells.Item(224, 2).Value2 = "china"
$ws.Cells.Item(224, 3).Value2 = "jia-league"
$ws.Cells.Item(224, 4).Value2 = "2023"
$ws.Cells.Item(224, 5).Value2 = 45227.35416666666
$ws.Cells.Item(224, 6).Value2 = "Nanjing City"
$ws.Cells.Item(224, 7).Value2 = 0
$ws.Cells.Item(224, 8).Value2 = "Sichuan Jiuniu"
$ws.Cells.Item(224, 9).Value2 = 1
$ws.Cells.Item(224, 10).Value2 = 2.76
$ws.Cells.Item(224, 11).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(224, 12).Value2 = 1.93
$ws.Cells.Item(224, 13).Value2 = "28/10/2023 08:27"
$ws.Cells.Item(224, 14).Value2 = 2.74
$ws.Cells.Item(224, 15).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(224, 16).Value2 = 3.18
$ws.Cells.Item(224, 17).Value2 = "28/10/2023 08:22"
$ws.Cells.Item(224, 18).Value2 = 2.68
$ws.Cells.Item(224, 19).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(224, 20).Value2 = 4.48
$ws.Cells.Item(224, 21).Value2 = "28/10/2023 08:27"
$ws.Cells.Item(224, 22).Value2 = "https://www.betexplorer.com/football/china/jia-league/nanjing-city-sichuan-jiuniu/hxQAnfaU/"

$ws.Cells.Item(225, 1).Value2 = 224
$ws.Cells.Item(225, 2).Value2 = "china"
$ws.Cells.Item(225, 3).Value2 = "jia-league"
$ws.Cells.Item(225, 4).Value2 = "2023"
$ws.Cells.Item(225, 5).Value2 = 45227.35416666666
$ws.Cells.Item(225, 6).Value2 = "Yanbian Longding"
$ws.Cells.Item(225, 7).Value2 = 3
$ws.Cells.Item(225, 8).Value2 = "Shanghai Jiading Huilong"
$ws.Cells.Item(225, 9).Value2 = 0
$ws.Cells.Item(225, 10).Value2 = 1.91
$ws.Cells.Item(225, 11).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(225, 12).Value2 = 1.29
$ws.Cells.Item(225, 13).Value2 = "28/10/2023 08:24"
$ws.Cells.Item(225, 14).Value2 = 3.15
$ws.Cells.Item(225, 15).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(225, 16).Value2 = 5.23
$ws.Cells.Item(225, 17).Value2 = "28/10/2023 08:25"
$ws.Cells.Item(225, 18).Value2 = 3.85
$ws.Cells.Item(225, 19).Value2 = "26/10/2023 20:42"
$ws.Cells.Item(225, 20).Value2 = 11.38
$ws.Cells.Item(225, 21).Value2 = "28/10/2023 08:25"
$ws.Cells.Item(225, 22).Value2 = "https://www.betexplorer.com/football/china/jia-league/yanbian-longding-shanghai-jiading-huilong/6Z2X3eU4/"

